$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 98935
$ws.Range("B3").Value = 98935
$ws.Range("B4").Value = 98935
$ws.Range("A5").Value = 130965940
$ws.Range("AC5").Value = 'Betydelsefulla förekomster . inventering åt vasa vind'
$ws.Range("AX5").Value = 'Pia Edfors, Anders Esplund, Enviro Planning'
$ws.Range("B5").Value = 98935
$ws.Range("Q5").Value = 496969
$ws.Range("R5").Value = 6713529
$ws.Range("A6").Value = 130965930
$ws.Range("AC6").Value = 'Måttlig förekomst . inventering åt vasa vind'
$ws.Range("AX6").Value = 'Pia Edfors, Enviro Planning'
$ws.Range("B6").Value = 98935
$ws.Range("Q6").Value = 496938
$ws.Range("R6").Value = 6713359
$ws.Range("B7").Value = 98935
$ws.Range("B8").Value = 98935
$ws.Range("B9").Value = 98935
$ws.Range("B10").Value = 98935
$ws.Range("A11").Value = 130965861
$ws.Range("AC11").Value = 'Betydelsefulla förekomster . inventering åt vasa vind'
$ws.Range("AX11").Value = 'Anders Esplund, Pia Edfors, Enviro Planning'
$ws.Range("B11").Value = 98935
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 219790
$ws.Range("F11").Value = 'Fläcknycklar'
$ws.Range("G11").Value = 'Dactylorhiza maculata'
$ws.Range("H11").Value = '(L.) Soó'
$ws.Range("Q11").Value = 497138
$ws.Range("R11").Value = 6713448
$ws.Range("A12").Value = 130965935
$ws.Range("AC12").Value = 'Måttlig förekomst . inventering åt vasa vind'
$ws.Range("AX12").Value = 'Pia Edfors, Enviro Planning'
$ws.Range("B12").Value = 79245
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("Q12").Value = 496969
$ws.Range("R12").Value = 6713674
$ws.Range("B13").Value = 98935
$ws.Range("B14").Value = 98935
